# Calificaciones Const y Var grupo B
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill grades (column E, "DE 0 A 5 CUMPLE?") for each student row.
# Row 4 (student 2) already had E4 = 4, leave it untouched.
$grades = @{
    3  = 5
    5  = 5
    6  = 5
    7  = 5
    8  = 5
    9  = 5
    10 = 5
    12 = 5
    13 = 5
    14 = 5
    15 = 5
    16 = 5
    17 = 5
    18 = 5
    19 = 5
    20 = 5
    21 = 5
    22 = 5
    23 = 5
    24 = 5
    25 = 5
    26 = 5
    28 = 5
}

foreach ($row in $grades.Keys) {
    $ws.Range("E$row").Value = $grades[$row]
}

# Row 27 (student 25) did not present: grade 0 and a remark in column F.
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = "No presenta"

# Update the active selection to match the saved workbook state.
$ws.Range("E4").Select()
